$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = "Datos actualizados a 11 de Abril de 2020 a las 14:22"

# Row 15: Paises Bajos
$ws.Cells.Item(15, 2).Value = 24413
$ws.Cells.Item(15, 3).Value = 1316
$ws.Cells.Item(15, 5).Value = 21520
$ws.Cells.Item(15, 7).Value = 132
$ws.Cells.Item(15, 8).Value = 2643

# Row 19: Austria
$ws.Cells.Item(19, 2).Value = 13767
$ws.Cells.Item(19, 3).Value = 207
$ws.Cells.Item(19, 5).Value = 6826

# Row 21: Israel
$ws.Cells.Item(21, 2).Value = 10525
$ws.Cells.Item(21, 3).Value = 117
$ws.Cells.Item(21, 4).Value = 1258
$ws.Cells.Item(21, 5).Value = 9171
$ws.Cells.Item(21, 6).Value = 180
$ws.Cells.Item(21, 7).Value = 1
$ws.Cells.Item(21, 8).Value = 96

# Row 23: Suecia
$ws.Cells.Item(23, 2).Value = 10151
$ws.Cells.Item(23, 3).Value = 466
$ws.Cells.Item(23, 5).Value = 8883
$ws.Cells.Item(23, 6).Value = 789
$ws.Cells.Item(23, 7).Value = 17
$ws.Cells.Item(23, 8).Value = 887

# Row 25: India
$ws.Cells.Item(25, 2).Value = 7876
$ws.Cells.Item(25, 3).Value = 276
$ws.Cells.Item(25, 5).Value = 6853

# Row 32: Dinamarca
$ws.Cells.Item(32, 1).Value = "Dinamarca"
$ws.Cells.Item(32, 2).Value = 5996
$ws.Cells.Item(32, 3).Value = 177
$ws.Cells.Item(32, 4).Value = 1955
$ws.Cells.Item(32, 5).Value = 3781
$ws.Cells.Item(32, 6).Value = 106
$ws.Cells.Item(32, 7).Value = 13
$ws.Cells.Item(32, 8).Value = 260

# Row 33: Rumania
$ws.Cells.Item(33, 1).Value = "Rumania"
$ws.Cells.Item(33, 2).Value = 5990
$ws.Cells.Item(33, 3).Value = 523
$ws.Cells.Item(33, 4).Value = 758
$ws.Cells.Item(33, 5).Value = 4950
$ws.Cells.Item(33, 6).Value = 208
$ws.Cells.Item(33, 7).Value = 12
$ws.Cells.Item(33, 8).Value = 282

# Row 34: Peru
$ws.Cells.Item(34, 1).Value = "Peru"
$ws.Cells.Item(34, 2).Value = 5897
$ws.Cells.Item(34, 4).Value = 1569
$ws.Cells.Item(34, 5).Value = 4159
$ws.Cells.Item(34, 6).Value = 130
$ws.Cells.Item(34, 8).Value = 169

# Row 60: Croacia
$ws.Cells.Item(60, 1).Value = "Croacia"
$ws.Cells.Item(60, 2).Value = 1534
$ws.Cells.Item(60, 3).Value = 39
$ws.Cells.Item(60, 4).Value = 323
$ws.Cells.Item(60, 5).Value = 1190
$ws.Cells.Item(60, 6).Value = 32
$ws.Cells.Item(60, 7).Value = 0
$ws.Cells.Item(60, 8).Value = 21

# Row 61: Marruecos
$ws.Cells.Item(61, 1).Value = "Marruecos"
$ws.Cells.Item(61, 2).Value = 1527
$ws.Cells.Item(61, 3).Value = 79
$ws.Cells.Item(61, 4).Value = 141
$ws.Cells.Item(61, 5).Value = 1276
$ws.Cells.Item(61, 6).Value = 1
$ws.Cells.Item(61, 7).Value = 3
$ws.Cells.Item(61, 8).Value = 110

# Row 77: Republica de Macedonia
$ws.Cells.Item(77, 1).Value = "Republica de Macedonia"
$ws.Cells.Item(77, 2).Value = 760
$ws.Cells.Item(77, 3).Value = 49
$ws.Cells.Item(77, 4).Value = 41
$ws.Cells.Item(77, 5).Value = 687
$ws.Cells.Item(77, 6).Value = 15
$ws.Cells.Item(77, 8).Value = 32

# Row 78: Uzbekistan
$ws.Cells.Item(78, 1).Value = "Uzbekistan"
$ws.Cells.Item(78, 2).Value = 729
$ws.Cells.Item(78, 3).Value = 105
$ws.Cells.Item(78, 4).Value = 42
$ws.Cells.Item(78, 5).Value = 684
$ws.Cells.Item(78, 6).Value = 8
$ws.Cells.Item(78, 8).Value = 3

# Row 79: Eslovaquia
$ws.Cells.Item(79, 1).Value = "Eslovaquia"
$ws.Cells.Item(79, 2).Value = 728
$ws.Cells.Item(79, 3).Value = 13
$ws.Cells.Item(79, 4).Value = 23
$ws.Cells.Item(79, 5).Value = 703
$ws.Cells.Item(79, 6).Value = 5
$ws.Cells.Item(79, 8).Value = 2

# Row 80: Crucero
$ws.Cells.Item(80, 1).Value = "Crucero"
$ws.Cells.Item(80, 2).Value = 712
$ws.Cells.Item(80, 4).Value = 619
$ws.Cells.Item(80, 5).Value = 82
$ws.Cells.Item(80, 6).Value = 10
$ws.Cells.Item(80, 8).Value = 11

# Row 107: Senegal
$ws.Cells.Item(107, 4).Value = 152
$ws.Cells.Item(107, 5).Value = 124

# Row 119: Republica de Yibuti
$ws.Cells.Item(119, 1).Value = "Republica de Yibuti"
$ws.Cells.Item(119, 2).Value = 187
$ws.Cells.Item(119, 3).Value = 37
$ws.Cells.Item(119, 4).Value = 36
$ws.Cells.Item(119, 5).Value = 149
$ws.Cells.Item(119, 7).Value = 1
$ws.Cells.Item(119, 8).Value = 2

# Row 120: Islas Feroe
$ws.Cells.Item(120, 1).Value = "Islas Feroe"
$ws.Cells.Item(120, 2).Value = 184
$ws.Cells.Item(120, 4).Value = 145
$ws.Cells.Item(120, 5).Value = 39
$ws.Cells.Item(120, 6).Value = 0
$ws.Cells.Item(120, 8).Value = 0

# Row 121: Venezuela
$ws.Cells.Item(121, 1).Value = "Venezuela"
$ws.Cells.Item(121, 2).Value = 175
$ws.Cells.Item(121, 4).Value = 84
$ws.Cells.Item(121, 5).Value = 82
$ws.Cells.Item(121, 6).Value = 6
$ws.Cells.Item(121, 8).Value = 9

# Row 122: Martinica
$ws.Cells.Item(122, 1).Value = "Martinica"
$ws.Cells.Item(122, 2).Value = 155
$ws.Cells.Item(122, 4).Value = 50
$ws.Cells.Item(122, 5).Value = 99
$ws.Cells.Item(122, 6).Value = 19
$ws.Cells.Item(122, 8).Value = 6
